$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.152.09"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "3.461.59"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'579.61"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Value = "'149.43"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'7.86"
$ws.Range("E9").Value = "  +2.53%  "
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").Value = "4.052.12"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").Value = "'28.56"
$ws.Range("E14").Value = "  -4.21%  "
$ws.Range("D15").Value = "3.450.25"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").Value = "63.164.57"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "'6.47"
$ws.Range("E18").Value = "  +2.96%  "
$ws.Range("D19").Value = "'14.41"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").Value = "'389.72"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").Value = "'0.562"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").Value = "'74.79"
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "3.593.67"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("E26").Value = "  -3.60%  "
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("D28").Value = "'7.68"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").Value = "'23.41"
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("E34").Value = "  -5.43%  "
$ws.Range("E35").Value = "  +3.15%  "
$ws.Range("D36").Value = "'5.34"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'32.09"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").Value = "'170.15"
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").Value = "3.499.08"
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("D41").Value = "'0.0774"
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("D42").Value = "'0.795"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").Value = "'42.85"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("E46").Value = "  -3.01%  "
$ws.Range("D47").Value = "2.585.94"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "'6.90"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").Value = "'22.67"
$ws.Range("E50").Value = "  -4.89%  "
$ws.Range("E51").Value = "  -0.04%  "
